$wb = $excel.ActiveWorkbook

# The PREGRADO sheet (first sheet) had an extra column L that held a
# "RATIFICAR EL ESTADO DE LA CARRERA..." ratification flag (with a cell
# comment on L5 and a data-validation dropdown on L6:L11). That column,
# along with its comment/validation/font, is removed in this revision.
$ws = $wb.Worksheets.Item("PREGRADO")
$ws.Columns.Item(12).Delete()
